$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Map of cell address -> new value
$changes = @{
    "E2"  = 65
    "E18" = 137
    "F18" = 63
    "H18" = 100
    "F20" = 4
    "H20" = 7
    "E26" = 37
    "F26" = 21
    "H26" = 31
    "E36" = 126
    "F36" = 61
    "H36" = 93
    "E37" = 63
    "F37" = 39
    "H37" = 51
    "E38" = 89
    "F38" = 22
    "H38" = 42
    "E40" = 29
    "E42" = 43
    "E44" = 34
    "F44" = 18
    "H44" = 28
    "E47" = 68
    "F47" = 43
    "H47" = 53
    "E48" = 43
    "F48" = 27
    "H48" = 33
    "E49" = 83
    "F49" = 45
    "H49" = 62
    "E55" = 10
    "E57" = 19
    "E62" = 56
    "E63" = 49
    "F63" = 20
    "H63" = 28
    "E64" = 39
    "F64" = 21
    "H64" = 26
    "E65" = 40
    "E70" = 51
    "E71" = 48
    "E72" = 53
    "F72" = 30
    "H72" = 41
    "E74" = 21
    "E76" = 60
    "E77" = 66
    "E79" = 48
    "F79" = 24
    "H79" = 35
    "E83" = 13
    "E84" = 7
    "E89" = 49
    "F89" = 24
    "H89" = 30
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}

$wb.Save()
